$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor reading row appended below the existing data (row 5).
$row = 5

# Force the numeric-looking values ("23", "0.0") to be stored as text,
# matching the existing rows where every column is textual data coming
# from the Adafruit IO feed export. Temporarily apply a text number
# format so Excel doesn't coerce the string to a number, then clear the
# formatting again so the new cells don't carry a stray style.
$dataRange = $ws.Range("C" + $row + ":F" + $row)
$dataRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T17:57:19Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "23"
$ws.Cells.Item($row, 4).Value = "0.0"
$ws.Cells.Item($row, 5).Value = "0.0"
$ws.Cells.Item($row, 6).Value = "0.0"

$dataRange.ClearFormats()
